# Auto-generated Excel COM-interop script to refresh market-price data cells
# on the Leve profit sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR), matching the
# scheduled-runner data refresh described in the commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1003
$arr[0,1] = 1003.3333
$arr[0,2] = 1002.5
$arr[0,3] = 1003.3333
$arr[0,4] = 1002.5
$ws.Range("H28:L28").Value = $arr
$ws.Range("M28").Value = -518.3333

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2701.9614
$arr[0,1] = 2695.5652
$arr[0,2] = 2751
$arr[0,3] = 2695.5652
$arr[0,4] = 2751
$ws.Range("H98:L98").Value = $arr
$ws.Range("M98").Value = -1197.5652

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1669.4445
$arr[0,1] = 624.125
$arr[0,2] = 10032
$arr[0,3] = 1872.375
$arr[0,4] = 30096
$ws.Range("H111:L111").Value = $arr
$ws.Range("M111").Value = 1194.625

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 3728.5789
$arr[0,1] = 3396
$arr[0,2] = 3747.0557
$arr[0,3] = 10188
$arr[0,4] = 11241.1671
$ws.Range("H112:L112").Value = $arr
$ws.Range("M112").Value = -9080
$ws.Range("N112").Value = -13457.1671

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2701.9614
$arr[0,1] = 2695.5652
$arr[0,2] = 2751
$arr[0,3] = 8086.6956
$arr[0,4] = 8253
$ws.Range("H122:L122").Value = $arr
$ws.Range("M122").Value = -5636.6956

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1597.875
$arr[0,1] = 663
$arr[0,2] = 1909.5
$arr[0,3] = 5967
$arr[0,4] = 17185.5
$ws.Range("H125:L125").Value = $arr
$ws.Range("M125").Value = -3507
$ws.Range("N125").Value = -22105.5

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2196.31
$arr[0,1] = 677.4
$arr[0,2] = 2365.078
$arr[0,3] = 2032.2
$arr[0,4] = 7095.234
$ws.Range("H129:L129").Value = $arr
$ws.Range("M129").Value = 2967.8
$ws.Range("N129").Value = -17095.234

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1980
$arr[0,1] = 1968.2
$arr[0,2] = 1999.6666
$arr[0,3] = 5904.6
$arr[0,4] = 5998.9998
$ws.Range("H131:L131").Value = $arr
$ws.Range("M131").Value = -864.6000000000004

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2921.1538
$arr[0,1] = 2915.9092
$arr[0,2] = 2950
$arr[0,3] = 8747.7276
$arr[0,4] = 8850
$ws.Range("H132:L132").Value = $arr
$ws.Range("M132").Value = -6217.7276

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 861
$arr[0,1] = 971.2857
$arr[0,2] = 475
$arr[0,3] = 8741.5713
$arr[0,4] = 4275
$ws.Range("H135:L135").Value = $arr
$ws.Range("M135").Value = -6206.5713

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2155.25
$arr[0,1] = 1231.6666
$arr[0,2] = 2709.4
$arr[0,3] = 3694.9998
$arr[0,4] = 8128.200000000001
$ws.Range("H137:L137").Value = $arr
$ws.Range("M137").Value = -1144.9998

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2776.4187
$arr[0,1] = 1511.5
$arr[0,2] = 3065.543
$arr[0,3] = 4534.5
$arr[0,4] = 9196.629000000001
$ws.Range("H138:L138").Value = $arr
$ws.Range("M138").Value = 605.5
$ws.Range("N138").Value = -19476.629

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2557347.8
$arr[0,1] = 2709890
$arr[0,2] = 675995.3
$arr[0,3] = 2709890
$arr[0,4] = 675995.3
$ws.Range("H32:L32").Value = $arr
$ws.Range("M32").Value = -2709603
$ws.Range("N32").Value = -676569.3

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 24592.182
$arr[0,1] = 23498
$arr[0,2] = 24835.334
$arr[0,3] = 23498
$arr[0,4] = 24835.334
$ws.Range("H43:L43").Value = $arr
$ws.Range("M43").Value = -23185
$ws.Range("N43").Value = -25461.334

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 3584.75
$arr[0,1] = 3336.4
$arr[0,2] = 3998.6667
$arr[0,3] = 3336.4
$arr[0,4] = 3998.6667
$ws.Range("H61:L61").Value = $arr
$ws.Range("M61").Value = -3124.4

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 20367.354
$arr[0,1] = 22719.732
$arr[0,2] = 2724.5
$arr[0,3] = 22719.732
$arr[0,4] = 2724.5
$ws.Range("H63:L63").Value = $arr
$ws.Range("M63").Value = -22033.732
$ws.Range("N63").Value = -4096.5

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 20367.354
$arr[0,1] = 22719.732
$arr[0,2] = 2724.5
$arr[0,3] = 113598.66
$arr[0,4] = 13622.5
$ws.Range("H66:L66").Value = $arr
$ws.Range("M66").Value = -110166.66
$ws.Range("N66").Value = -20486.5

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2536
$arr[0,1] = 2360.4
$arr[0,2] = 2770.1333
$arr[0,3] = 2360.4
$arr[0,4] = 2770.1333
$ws.Range("H74:L74").Value = $arr
$ws.Range("M74").Value = -1486.4

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2536
$arr[0,1] = 2360.4
$arr[0,2] = 2770.1333
$arr[0,3] = 11802
$arr[0,4] = 13850.6665
$ws.Range("H77:L77").Value = $arr
$ws.Range("M77").Value = -7434

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 44625.25
$arr[0,1] = 0
$arr[0,2] = 44625.25
$arr[0,3] = 0
$arr[0,4] = 44625.25
$ws.Range("H101:L101").Value = $arr
$ws.Range("N101").Value = -51115.25

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2289.158
$arr[0,1] = 2470.4614
$arr[0,2] = 1896.3334
$arr[0,3] = 7411.3842
$arr[0,4] = 5689.0002
$ws.Range("H132:L132").Value = $arr
$ws.Range("M132").Value = -4881.3842

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 3584.75
$arr[0,1] = 3336.4
$arr[0,2] = 3998.6667
$arr[0,3] = 10009.2
$arr[0,4] = 11996.0001
$ws.Range("H136:L136").Value = $arr
$ws.Range("M136").Value = -7459.200000000001

$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 206427.1
$arr[0,1] = 9176
$arr[0,2] = 403678.2
$arr[0,3] = 9176
$arr[0,4] = 403678.2
$ws.Range("H99:L99").Value = $arr
$ws.Range("M99").Value = -7678
$ws.Range("N99").Value = -406674.2

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 4419.263
$arr[0,1] = 4523.1113
$arr[0,2] = 2550
$arr[0,3] = 13569.3339
$arr[0,4] = 7650
$ws.Range("H134:L134").Value = $arr
$ws.Range("M134").Value = -11034.3339

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$ws.Range("H31:L31").Value = $arr
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$ws.Range("H34:L34").Value = $arr
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 3178.303
$arr[0,1] = 2693.75
$arr[0,2] = 3923.7693
$arr[0,3] = 2693.75
$arr[0,4] = 3923.7693
$ws.Range("H58:L58").Value = $arr
$ws.Range("M58").Value = -2490.75

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 30491
$arr[0,1] = 39125
$arr[0,2] = 24735
$arr[0,3] = 39125
$arr[0,4] = 24735
$ws.Range("H60:L60").Value = $arr
$ws.Range("M60").Value = -38614

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 46137.75
$arr[0,1] = 55000
$arr[0,2] = 43183.668
$arr[0,3] = 55000
$arr[0,4] = 43183.668
$ws.Range("H131:L131").Value = $arr
$ws.Range("M131").Value = -49960
$ws.Range("N131").Value = -53263.668

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2169.577
$arr[0,1] = 2196.36
$arr[0,2] = 1500
$arr[0,3] = 6589.08
$arr[0,4] = 4500
$ws.Range("H132:L132").Value = $arr
$ws.Range("M132").Value = -4059.08
$ws.Range("N132").Value = -9560

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1631.963
$arr[0,1] = 1651.0869
$arr[0,2] = 1522
$arr[0,3] = 4953.2607
$arr[0,4] = 4566
$ws.Range("H134:L134").Value = $arr
$ws.Range("M134").Value = -2418.2607

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 3178.303
$arr[0,1] = 2693.75
$arr[0,2] = 3923.7693
$arr[0,3] = 8081.25
$arr[0,4] = 11771.3079
$ws.Range("H136:L136").Value = $arr
$ws.Range("M136").Value = -5531.25

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 80539.16
$arr[0,1] = 1793
$arr[0,2] = 91277.27
$arr[0,3] = 5379
$arr[0,4] = 273831.81
$ws.Range("H5:L5").Value = $arr
$ws.Range("M5").Value = -5267
$ws.Range("N5").Value = -274055.81

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$ws.Range("H56:L56").Value = $arr
$ws.Range("M56").ClearContents()

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2721.25
$arr[0,1] = 1295
$arr[0,2] = 7000
$arr[0,3] = 3885
$arr[0,4] = 21000
$ws.Range("H57:L57").Value = $arr
$ws.Range("M57").Value = -3326
$ws.Range("N57").Value = -22118

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 14539.429
$arr[0,1] = 13000
$arr[0,2] = 14796
$arr[0,3] = 39000
$arr[0,4] = 44388
$ws.Range("H93:L93").Value = $arr
$ws.Range("M93").Value = -37128
$ws.Range("N93").Value = -48132

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 80539.16
$arr[0,1] = 1793
$arr[0,2] = 91277.27
$arr[0,3] = 16137
$arr[0,4] = 821495.4300000001
$ws.Range("H135:L135").Value = $arr
$ws.Range("M135").Value = -13602
$ws.Range("N135").Value = -826565.4300000001

$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 936.6667
$arr[0,1] = 966.9167
$arr[0,2] = 815.6667
$arr[0,3] = 966.9167
$arr[0,4] = 815.6667
$ws.Range("H102:L102").Value = $arr
$ws.Range("M102").Value = 655.0833
$ws.Range("N102").Value = -4059.6667

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1272.0454
$arr[0,1] = 1115.5264
$arr[0,2] = 2263.3333
$arr[0,3] = 3346.5792
$arr[0,4] = 6789.999899999999
$ws.Range("H132:L132").Value = $arr
$ws.Range("M132").Value = -816.5792000000001

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 489.25
$arr[0,1] = 485.66666
$arr[0,2] = 500
$arr[0,3] = 485.66666
$arr[0,4] = 500
$ws.Range("H55:L55").Value = $arr
$ws.Range("M55").Value = -312.66666
$ws.Range("N55").Value = -846

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 17048.111
$arr[0,1] = 0
$arr[0,2] = 17048.111
$arr[0,3] = 0
$arr[0,4] = 17048.111
$ws.Range("H106:L106").Value = $arr
$ws.Range("N106").Value = -19572.111

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 44999
$arr[0,1] = 0
$arr[0,2] = 44999
$arr[0,3] = 0
$arr[0,4] = 44999
$ws.Range("H129:L129").Value = $arr
$ws.Range("N129").Value = -54999

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 3540.5
$arr[0,1] = 2988.25
$arr[0,2] = 5749.5
$arr[0,3] = 8964.75
$arr[0,4] = 17248.5
$ws.Range("H132:L132").Value = $arr
$ws.Range("M132").Value = -6434.75

$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2918.125
$arr[0,1] = 3251.25
$arr[0,2] = 2807.0833
$arr[0,3] = 3251.25
$arr[0,4] = 2807.0833
$ws.Range("H96:L96").Value = $arr
$ws.Range("M96").Value = -1878.25
$ws.Range("N96").Value = -5553.0833

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1714.6666
$arr[0,1] = 1722.75
$arr[0,2] = 1650
$arr[0,3] = 5168.25
$arr[0,4] = 4950
$ws.Range("H126:L126").Value = $arr
$ws.Range("M126").Value = -2698.25
$ws.Range("N126").Value = -9890

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2124.8293
$arr[0,1] = 2252.0293
$arr[0,2] = 1507
$arr[0,3] = 6756.0879
$arr[0,4] = 4521
$ws.Range("H132:L132").Value = $arr
$ws.Range("M132").Value = -4226.0879
$ws.Range("N132").Value = -9581

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 830.3929000000001
$arr[0,1] = 830.3929000000001
$arr[0,2] = 0
$arr[0,3] = 2491.1787
$arr[0,4] = 0
$ws.Range("H136:L136").Value = $arr
$ws.Range("M136").Value = 58.82129999999961
